$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: new - ItemC / UnitC
$ws.Range("A5").Value = "ItemC"
$ws.Range("B5").Value = "UnitC"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 4

# Row 6: new - UnitD entered before ItemD
$ws.Range("B6").Value = "UnitD"
$ws.Range("A6").Value = "ItemD"
$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = 2

# Row 3: ItemA -> Item,A (note comma) -- edited last
$ws.Range("A3").Value = "Item,A"

# Selection on A3 (matches sheetView selection in diff)
$ws.Range("A3").Select()
